$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Etat Virement")

# --- Row 2 (NASIRI HASNAA) : only amounts change ---
$ws.Range("I2").Value = 9999.99
$ws.Range("J2").Value = 999.99
$ws.Range("K2").Value = 9000

# --- Row 3 : becomes ZERNAKH ABDELLAH / 052 AV1 / 3000-450-2550 ---
$ws.Range("A3").Value = "ZERNAKH ABDELLAH"
$ws.Range("B3").Value = "IB19558"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "145101211406073828000084"
$ws.Range("D3").Value = "MARRAKECH BENI MELLAL"
$ws.Range("E3").Value = "BP"
$ws.Range("F3").Value = "Point de vente"
$ws.Range("G3").Value = "052/FKIH BEN SALEH/AV1"
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 3000
$ws.Range("J3").Value = 450
$ws.Range("K3").Value = 2550

# --- Row 4 : becomes ZERNAKH ABDELLAH / 052 AV1 / 12000-0-12000 ---
$ws.Range("A4").Value = "ZERNAKH ABDELLAH"
$ws.Range("B4").Value = "IB19558"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "145101211406073828000084"
$ws.Range("D4").Value = "MARRAKECH BENI MELLAL"
$ws.Range("E4").Value = "BP"
$ws.Range("F4").Value = "Point de vente"
$ws.Range("G4").Value = "052/FKIH BEN SALEH/AV1"
$ws.Range("H4").Value = "mensuelle"
$ws.Range("I4").Value = 12000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 12000

# --- Row 5 : becomes MOHAMED BADRANE / KHOURIBGA / 7500-375-7125 ---
$ws.Range("A5").Value = "MOHAMED BADRANE"
$ws.Range("B5").Value = "I83603"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "225400000805987601012173"
$ws.Range("D5").Value = "KHOURIBGA"
$ws.Range("E5").Value = "CA"
$ws.Range("F5").Value = "Point de vente"
$ws.Range("G5").Value = "605/KHOURIBGA NAHDA"
$ws.Range("H5").Value = "mensuelle"
$ws.Range("I5").Value = 7500
$ws.Range("J5").Value = 375
$ws.Range("K5").Value = 7125

# --- Row 8 totals ---
$ws.Range("I8").Value = 45999.99
$ws.Range("J8").Value = 3174.99
$ws.Range("K8").Value = 42825
